$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds the "last changed" date for each record.
# Every data row (2 through 115) is being refreshed from 2023-09-09 (45178)
# to 2023-09-10 (45179).
$lastRow = 115
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45178) {
        $cell.Value = 45179
    }
}
